$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D9"  = -8.084399999999993
    "D13" = -8.063999999999998
    "D16" = -8.146099999999997
    "D18" = -8.551600000000002
    "D20" = -7.074800000000003
    "D26" = -7.663499999999999
    "D27" = -7.801599999999999
    "D29" = -7.318899999999998
    "D35" = -7.250200000000006
    "D36" = -7.752299999999998
    "D45" = -7.090299999999997
    "D55" = -8.826299999999998
    "D57" = -8.310699999999999
    "D69" = -7.521899999999996
    "D76" = -7.4691
    "D78" = -7.293800000000001
    "D82" = -8.368700000000004
    "D83" = -8.9658
    "D93" = -6.996399999999992
    "D97" = -7.693699999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
